# Add six new event rows (378-383) to the "Tabelle1" worksheet, each with a
# date, event name, location, city and a hyperlinked ticket/social-media URL.
# This mirrors the existing pattern used by all prior rows in the sheet
# (column A = date serial with style "dd.mm.yy", columns B-E = left aligned
# text cells, column E additionally carrying an external hyperlink).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellref, $text) {
    $r = $ws.Range($cellref)
    $r.NumberFormat = "@"
    $r.Value = $text
}

function Set-DateCell($cellref, $serial) {
    $r = $ws.Range($cellref)
    $r.Value = $serial
}

function Set-HyperlinkCell($cellref, $url) {
    $r = $ws.Range($cellref)
    $r.NumberFormat = "@"
    $r.Value = $url
    $ws.Hyperlinks.Add($r, $url, "", "", $url)
    # Hyperlinks.Add stamps the cell with a brand new "hyperlink" font
    # (theme colour + underline). Reset the font back to the sheet's normal
    # look so the cell keeps using the same plain text style as columns B-D,
    # just like in the rows that were already in the workbook.
    $f = $r.Font
    $f.Name = "Calibri"
    $f.Size = 11
    $f.Underline = 0
    $f.ColorIndex = 0
    $f.Color = 0
    $r.NumberFormat = "@"
}

# Row 378 - LIEBEFELD TANZ IN DEN MAI
Set-DateCell "A378" 45777
Set-TextCell "B378" "LIEBEFELD TANZ IN DEN MAI"
Set-TextCell "C378" "Sams"
Set-TextCell "D378" "Bielefeld"
Set-HyperlinkCell "E378" "https://www.instagram.com/reel/DH3A7HoqFwJ/?igsh=MTRxNTA5bGgwbXd0Nw=="

# Row 379 - NIGHTROOMS x SIXSIXSOUNDS
Set-DateCell "A379" 45786
Set-TextCell "B379" "NIGHTROOMS x SIXSIXSOUNDS"
Set-TextCell "C379" "Nighrooms"
Set-TextCell "D379" "Dortmund"
Set-HyperlinkCell "E379" "https://www.instagram.com/reel/DIekUsbsCHK/?igsh=MTBiNHhjODJvcmZzeg=="

# Row 380 - DESIRE
Set-DateCell "A380" 45781
Set-TextCell "B380" "DESIRE"
Set-TextCell "C380" "Junkyard"
Set-TextCell "D380" "Dortmund"
Set-HyperlinkCell "E380" "https://www.instagram.com/reel/DH80zwBsIrU/?igsh=bzBrYXlrM28xcnNx"

# Row 381 - TECHNO DAMPFER DUISBURG
Set-DateCell "A381" 45777
Set-TextCell "B381" "TECHNO DAMPFER DUISBURG"
Set-TextCell "C381" "Mercator Insel"
Set-TextCell "D381" "Duisburg"
Set-HyperlinkCell "E381" "https://technodampfer.ticket.io/?view=table"

# Row 382 - TECHNO DAMPFER DÜSSELDORF
Set-DateCell "A382" 45778
Set-TextCell "B382" "TECHNO DAMPFER DÜSSELDORF"
Set-TextCell "C382" "Tonhallenufer"
Set-TextCell "D382" "Düsseldorf"
Set-HyperlinkCell "E382" "https://technodampfer.ticket.io/?view=table"

# Row 383 - TECHNO DAMPFER KÖLN
Set-DateCell "A383" 45779
Set-TextCell "B383" "TECHNO DAMPFER KÖLN"
Set-TextCell "C383" "Am Leystapel"
Set-TextCell "D383" "Köln"
Set-HyperlinkCell "E383" "https://technodampfer.ticket.io/?view=table"
